# 21-02-2023 after mdec meeting
# Add "ic number" (new column D), " bank name" and " bank account number"
# (new trailing columns) to the Employee List import template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee List")

# 1. Append the two bank columns at the current last position (P, Q) first,
#    so the shared-string table gets the same allocation order as the
#    original edit (bank name / bank account number before ic number).
$ws.Range("O2").Copy()
$ws.Range("P2:Q2").PasteSpecial(-4122)
$ws.Range("P2").Value = " bank name"
$ws.Range("Q2").Value = " bank account number"

# 2. Insert a new column before D (shifts D:Q -> E:R) for the "ic number"
#    field, right after the employee's name.
$ws.Columns.Item(4).Insert()
$ws.Range("D2").Value = "ic number"

# 3. The "country" header (with its hyperlink to the reference sheet) moved
#    from G2 to H2 because of the column insert above. The stored hyperlink
#    anchor does not shift automatically, so re-create it at H2. Re-adding a
#    hyperlink resets the cell look, so stash/restore the original
#    formatting (white text on the dark header fill) around the call.
$ws.Range("H2").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("G2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "", "'References for Country Name'!A1", "", "country")
$ws.Range("Z100").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# 4. Extend the title merge across the two new trailing columns.
$ws.Range("A1:P1").UnMerge()
$ws.Range("A1:R1").Merge()

# 5. Mirror the blank input-row formatting from column B into the new
#    column C for the sample rows.
$ws.Range("B3").Copy()
$ws.Range("C3:C5").PasteSpecial(-4122)

# 6. Match the saved view state (scrolled right, new selected cell).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P15").Select()
